$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (columns A..F) after reordering
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "living_rooms_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "bedrooms_2"
$ws.Range("E1").Value = "kitchens_1"
$ws.Range("F1").Value = "kitchens_2"

# Updated one-hot data rows 2..7 (columns A..F)
$data = @(
    @(0,0,0,1,0,0),
    @(0,0,0,0,1,0),
    @(1,0,0,0,0,0),
    @(0,0,1,0,0,0),
    @(0,0,0,0,0,1),
    @(0,1,0,0,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
